$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value2 = 0
$ws.Range("I43").Value2 = 0
$ws.Range("J43").Value2 = 0
$ws.Range("K43").Value2 = 0
$ws.Range("L43").Value2 = 0
$ws.Range("H132").Value2 = 16677.578
$ws.Range("I132").Value2 = 3034
$ws.Range("J132").Value2 = 60640.223
$ws.Range("K132").Value2 = 9102
$ws.Range("L132").Value2 = 181920.669
$ws.Range("M132").Value2 = -6572
$ws.Range("N132").Value2 = -186980.669
$ws.Range("H135").Value2 = 16130620
$ws.Range("I135").Value2 = 17858472
$ws.Range("K135").Value2 = 160726248
$ws.Range("M135").Value2 = -160723713
$ws.Range("H137").Value2 = 3076.0977
$ws.Range("I137").Value2 = 2564.6858
$ws.Range("K137").Value2 = 7694.057400000001
$ws.Range("M137").Value2 = -5144.057400000001
$ws.Range("H138").Value2 = 2723.5073
$ws.Range("I138").Value2 = 1353.5883
$ws.Range("J138").Value2 = 4054.2856
$ws.Range("K138").Value2 = 4060.7649
$ws.Range("L138").Value2 = 12162.8568
$ws.Range("M138").Value2 = 1079.2351
$ws.Range("N138").Value2 = -22442.8568
$ws.Range("H141").Value2 = 3466.2144
$ws.Range("I141").Value2 = 2259.32
$ws.Range("K141").Value2 = 6777.960000000001
$ws.Range("M141").Value2 = -1597.960000000001
$ws.Range("M43").ClearContents() | Out-Null
$ws.Range("N43").ClearContents() | Out-Null

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 20835176
$ws.Range("I32").Value2 = 22224096
$ws.Range("K32").Value2 = 22224096
$ws.Range("M32").Value2 = -22223809
$ws.Range("H61").Value2 = 792.94116
$ws.Range("I61").Value2 = 786.6667
$ws.Range("J61").Value2 = 1000
$ws.Range("K61").Value2 = 786.6667
$ws.Range("L61").Value2 = 1000
$ws.Range("M61").Value2 = -574.6667
$ws.Range("N61").Value2 = -1424
$ws.Range("H132").Value2 = 2564.6667
$ws.Range("I132").Value2 = 2564.6667
$ws.Range("K132").Value2 = 7694.000100000001
$ws.Range("M132").Value2 = -5164.000100000001
$ws.Range("H136").Value2 = 792.94116
$ws.Range("I136").Value2 = 786.6667
$ws.Range("J136").Value2 = 1000
$ws.Range("K136").Value2 = 2360.0001
$ws.Range("L136").Value2 = 3000
$ws.Range("M136").Value2 = 189.9998999999998
$ws.Range("N136").Value2 = -8100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 1461.4667
$ws.Range("I20").Value2 = 1409.5
$ws.Range("J20").Value2 = 1565.4
$ws.Range("K20").Value2 = 1409.5
$ws.Range("L20").Value2 = 1565.4
$ws.Range("M20").Value2 = -1162.5
$ws.Range("N20").Value2 = -2059.4
$ws.Range("H94").Value2 = 2173.7407
$ws.Range("I94").Value2 = 2591.111
$ws.Range("K94").Value2 = 2591.111
$ws.Range("M94").Value2 = -2140.111
$ws.Range("H105").Value2 = 1661.6
$ws.Range("I105").Value2 = 1661.6
$ws.Range("J105").Value2 = 0
$ws.Range("K105").Value2 = 1661.6
$ws.Range("L105").Value2 = 0
$ws.Range("M105").Value2 = 85.40000000000009
$ws.Range("H134").Value2 = 1398.3235
$ws.Range("I134").Value2 = 1228.5758
$ws.Range("K134").Value2 = 3685.7274
$ws.Range("M134").Value2 = -1150.7274
$ws.Range("H138").Value2 = 72147.92
$ws.Range("J138").Value2 = 72147.92
$ws.Range("L138").Value2 = 72147.92
$ws.Range("N138").Value2 = -82427.92
$ws.Range("N105").ClearContents() | Out-Null

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1917.1305
$ws.Range("I31").Value2 = 1822.4546
$ws.Range("K31").Value2 = 1822.4546
$ws.Range("M31").Value2 = -1527.4546
$ws.Range("H34").Value2 = 1917.1305
$ws.Range("I34").Value2 = 1822.4546
$ws.Range("K34").Value2 = 1822.4546
$ws.Range("M34").Value2 = -1620.4546
$ws.Range("H132").Value2 = 2436.8333
$ws.Range("J132").Value2 = 3170.8333
$ws.Range("L132").Value2 = 9512.499899999999
$ws.Range("N132").Value2 = -14572.4999
$ws.Range("H135").Value2 = 70000
$ws.Range("J135").Value2 = 70000
$ws.Range("L135").Value2 = 70000
$ws.Range("N135").Value2 = -80140

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value2 = 2077.3
$ws.Range("I137").Value2 = 1863.8889
$ws.Range("J137").Value2 = 3998
$ws.Range("K137").Value2 = 5591.6667
$ws.Range("L137").Value2 = 11994
$ws.Range("M137").Value2 = -491.6666999999998
$ws.Range("N137").Value2 = -22194

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 142331.12
$ws.Range("J70").Value2 = 3666.6667
$ws.Range("L70").Value2 = 3666.6667
$ws.Range("H73").Value2 = 142331.12
$ws.Range("J73").Value2 = 3666.6667
$ws.Range("L73").Value2 = 3666.6667
$ws.Range("H102").Value2 = 2937.375
$ws.Range("I102").Value2 = 2010.7
$ws.Range("K102").Value2 = 2010.7
$ws.Range("M102").Value2 = -388.7
$ws.Range("H117").Value2 = 49500
$ws.Range("J117").Value2 = 49500
$ws.Range("L117").Value2 = 49500
$ws.Range("N117").Value2 = -56384
$ws.Range("H132").Value2 = 1966.6666
$ws.Range("I132").Value2 = 1865
$ws.Range("K132").Value2 = 5595
$ws.Range("M132").Value2 = -3065
$ws.Range("N70").Value2 = -4206.6667
$ws.Range("N73").Value2 = -5538.6667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 600
$ws.Range("J46").Value2 = 500
$ws.Range("L46").Value2 = 500
$ws.Range("N46").Value2 = -876
$ws.Range("H61").Value2 = 1690.6428
$ws.Range("I61").Value2 = 1772.5
$ws.Range("K61").Value2 = 1772.5
$ws.Range("M61").Value2 = -1570.5
$ws.Range("H113").Value2 = 1690.6428
$ws.Range("I113").Value2 = 1772.5
$ws.Range("K113").Value2 = 1772.5
$ws.Range("M113").Value2 = 397.5
$ws.Range("H118").Value2 = 50000
$ws.Range("J118").Value2 = 50000
$ws.Range("L118").Value2 = 50000
$ws.Range("N118").Value2 = -53314
$ws.Range("H132").Value2 = 2599.9302
$ws.Range("I132").Value2 = 2303.2942
$ws.Range("J132").Value2 = 3720.5557
$ws.Range("K132").Value2 = 6909.882599999999
$ws.Range("L132").Value2 = 11161.6671
$ws.Range("M132").Value2 = -4379.882599999999
$ws.Range("N132").Value2 = -16221.6671
$ws.Range("H136").Value2 = 2767
$ws.Range("J136").Value2 = 3206.1428
$ws.Range("L136").Value2 = 9618.428400000001
$ws.Range("N136").Value2 = -14718.4284

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value2 = 7568.3335
$ws.Range("I20").Value2 = 14705
$ws.Range("K20").Value2 = 14705
$ws.Range("M20").Value2 = -14465
$ws.Range("H62").Value2 = 3587.8
$ws.Range("I62").Value2 = 3649
$ws.Range("J62").Value2 = 3572.5
$ws.Range("K62").Value2 = 3649
$ws.Range("L62").Value2 = 3572.5
$ws.Range("M62").Value2 = -3025
$ws.Range("N62").Value2 = -4820.5
$ws.Range("H65").Value2 = 3587.8
$ws.Range("I65").Value2 = 3649
$ws.Range("J65").Value2 = 3572.5
$ws.Range("K65").Value2 = 18245
$ws.Range("L65").Value2 = 17862.5
$ws.Range("M65").Value2 = -15125
$ws.Range("N65").Value2 = -24102.5
$ws.Range("H81").Value2 = 15877391
$ws.Range("I81").Value2 = 1875
$ws.Range("K81").Value2 = 3750
$ws.Range("M81").Value2 = -2689
$ws.Range("H84").Value2 = 15877391
$ws.Range("I84").Value2 = 1875
$ws.Range("K84").Value2 = 18750
$ws.Range("M84").Value2 = -13446
$ws.Range("H126").Value2 = 1709.375
$ws.Range("I126").Value2 = 1605.8182
$ws.Range("J126").Value2 = 1937.2
$ws.Range("K126").Value2 = 4817.4546
$ws.Range("L126").Value2 = 5811.6
$ws.Range("M126").Value2 = -2347.4546
$ws.Range("N126").Value2 = -10751.6
$ws.Range("H132").Value2 = 1481.3529
$ws.Range("I132").Value2 = 1499.6875
$ws.Range("K132").Value2 = 4499.0625
$ws.Range("M132").Value2 = -1969.0625
$ws.Range("H135").Value2 = 67993.836
$ws.Range("J135").Value2 = 67993.836
$ws.Range("L135").Value2 = 67993.836
$ws.Range("N135").Value2 = -78133.836
